# Generate Report for Handoff
# Updates the Overview / zh-cn / de-de sheets with the latest handoff-report
# data: two previously-handed-off files move to "Ready for handoff" status
# (one .md renamed/replaced, one .md -> new uuid.png), and a brand new
# cf007469-....png file is appended as row 4 on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Clear any pre-existing hyperlinks on the sheet so we can rebuild them
# (this engine's Hyperlinks.Delete() clears the whole sheet's collection).
$ov.Range("A2").Hyperlinks.Delete()

# Row 2
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-03-23 10:08:26"
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/8c79659b-3b19-4976-a1b7-2471bcce5cee.png", [Type]::Missing, [Type]::Missing, "8c79659b-3b19-4976-a1b7-2471bcce5cee.png") | Out-Null

# Row 3
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-23 10:08:26"
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/cd364228-21a2-4e1f-80e9-750c7160fcba.md", [Type]::Missing, [Type]::Missing, "cd364228-21a2-4e1f-80e9-750c7160fcba.md") | Out-Null

# Row 4 (new)
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-23 10:08:26"
$ov.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/cf007469-bcf1-4ec9-9f4f-9d0e3e7c2940.png", [Type]::Missing, [Type]::Missing, "cf007469-bcf1-4ec9-9f4f-9d0e3e7c2940.png") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Hyperlinks.Delete()

# Row 2: source file is now the .png (Ready for handoff, no Latest Target
# File / Latest Handback File / Reference Tokens any more -> drop F,G,I,L)
$zh.Range("B2").Value = ".png"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("F2").Clear()
$zh.Range("G2").Clear()
$zh.Range("E2").Value = "2016-03-23 10:08:18"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "IsDependency"
$zh.Range("K2").Value = "e2e\cd364228-21a2-4e1f-80e9-750c7160fcba.md"
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/8c79659b-3b19-4976-a1b7-2471bcce5cee.png", [Type]::Missing, [Type]::Missing, "8c79659b-3b19-4976-a1b7-2471bcce5cee.png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5e1f3548be4c751df4da367aa64b68136fd0859/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/8d82dfa52fe1f50b9263de0baa0920e81b619f27.png", [Type]::Missing, [Type]::Missing, "8d82dfa52fe1f50b9263de0baa0920e81b619f27.png") | Out-Null

# Row 3
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Clear()
$zh.Range("G3").Clear()
$zh.Range("E3").Value = "2016-03-23 10:08:18"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/cd364228-21a2-4e1f-80e9-750c7160fcba.md", [Type]::Missing, [Type]::Missing, "cd364228-21a2-4e1f-80e9-750c7160fcba.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5e1f3548be4c751df4da367aa64b68136fd0859/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/cd364228-21a2-4e1f-80e9-750c7160fcba.9b7767f8ecd045ebf40cf7a3de6be87b630813d6.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "cd364228-21a2-4e1f-80e9-750c7160fcba.9b7767f8ecd045ebf40cf7a3de6be87b630813d6.zh-cn.xlf") | Out-Null

# Row 4 (new)
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("E4").Value = "2016-03-23 10:08:18"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("J4").Value = "IsDependency"
$zh.Range("K4").Value = "e2e\cd364228-21a2-4e1f-80e9-750c7160fcba.md"
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/cf007469-bcf1-4ec9-9f4f-9d0e3e7c2940.png", [Type]::Missing, [Type]::Missing, "cf007469-bcf1-4ec9-9f4f-9d0e3e7c2940.png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5e1f3548be4c751df4da367aa64b68136fd0859/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/b420084c1625a0910ef6b00d2bae93e6764b5b6a.png", [Type]::Missing, [Type]::Missing, "b420084c1625a0910ef6b00d2bae93e6764b5b6a.png") | Out-Null

$zh.Range("E2:E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H2:H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Hyperlinks.Delete()

# Row 2
$de.Range("B2").Value = ".png"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("F2").Clear()
$de.Range("G2").Clear()
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "IsDependency"
$de.Range("K2").Value = "e2e\cd364228-21a2-4e1f-80e9-750c7160fcba.md"
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/8c79659b-3b19-4976-a1b7-2471bcce5cee.png", [Type]::Missing, [Type]::Missing, "8c79659b-3b19-4976-a1b7-2471bcce5cee.png") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a20dcbfe7cabc7e8393b2b59f182b112255dd4b8/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/8d82dfa52fe1f50b9263de0baa0920e81b619f27.png", [Type]::Missing, [Type]::Missing, "8d82dfa52fe1f50b9263de0baa0920e81b619f27.png") | Out-Null

# Row 3
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Clear()
$de.Range("G3").Clear()
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/cd364228-21a2-4e1f-80e9-750c7160fcba.md", [Type]::Missing, [Type]::Missing, "cd364228-21a2-4e1f-80e9-750c7160fcba.md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a20dcbfe7cabc7e8393b2b59f182b112255dd4b8/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/cd364228-21a2-4e1f-80e9-750c7160fcba.9b7767f8ecd045ebf40cf7a3de6be87b630813d6.de-de.xlf", [Type]::Missing, [Type]::Missing, "cd364228-21a2-4e1f-80e9-750c7160fcba.9b7767f8ecd045ebf40cf7a3de6be87b630813d6.de-de.xlf") | Out-Null

# Row 4 (new)
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("E4").Value = "2016-03-23 10:08:26"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("J4").Value = "IsDependency"
$de.Range("K4").Value = "e2e\cd364228-21a2-4e1f-80e9-750c7160fcba.md"
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/53e4d6a04be53eee20341544517718a177087b4e/e2e/cf007469-bcf1-4ec9-9f4f-9d0e3e7c2940.png", [Type]::Missing, [Type]::Missing, "cf007469-bcf1-4ec9-9f4f-9d0e3e7c2940.png") | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a20dcbfe7cabc7e8393b2b59f182b112255dd4b8/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/b420084c1625a0910ef6b00d2bae93e6764b5b6a.png", [Type]::Missing, [Type]::Missing, "b420084c1625a0910ef6b00d2bae93e6764b5b6a.png") | Out-Null

$de.Range("H2:H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
